# Logged Week 15 and simulated Week 16
# Appends this week's per-play yardage logs to the YDS and ST running
# tally cells, and updates the season-to-date totals on OFF / DEF / ST /
# TURNS / PEN accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: append a space-separated list of numbers to the existing
# space-separated text already stored in a cell.
# ---------------------------------------------------------------------
function Append-Numbers($ws, [string]$cellRef, [string]$numbersToAdd) {
    $existing = $ws.Range($cellRef).Value2
    $ws.Range($cellRef).Value = ($existing + " " + $numbersToAdd)
}

# ---------------------------------------------------------------------
# YDS sheet: append Week 16 per-play yard logs
# ---------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

Append-Numbers $ydsWs "B2" "4 -7 -2 8 1 4 4 -1 2 4 24 2 4 1 1 6 4 3 12 2 6 2 5 0"
Append-Numbers $ydsWs "C2" "2 -2 1 3 19 7 6 3 3 6 3 5 8 6 3 3 0 1 6 2 9 4 1 -3"
Append-Numbers $ydsWs "B3" "5 10 6 6 2 4 3 4 4 9 13 6 6 13 5 25 11 3 6"
Append-Numbers $ydsWs "C3" "-1 12 14 11 17 12 7 5 5 10 5 13 8 6 16 14 16 12 6 7 3 12 9 12 15"

# ---------------------------------------------------------------------
# OFF sheet: updated season totals
# ---------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("C2").Value = 187
$offWs.Range("E2").Value = 11
$offWs.Range("F2").Value = 61
$offWs.Range("H2").Value = 6
$offWs.Range("J2").Value = 29
$offWs.Range("L2").Value = 237
$offWs.Range("M2").Value = 156
$offWs.Range("O2").Value = 24
$offWs.Range("P2").Value = 10
$offWs.Range("Q2").Value = 495

$offWs.Range("B3").Value = 14
$offWs.Range("C3").Value = 124
$offWs.Range("F3").Value = 84
$offWs.Range("G3").Value = 24
$offWs.Range("H3").Value = 27
$offWs.Range("I3").Value = 42
$offWs.Range("J3").Value = 37

# ---------------------------------------------------------------------
# DEF sheet: updated season totals
# ---------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("B2").Value = 7
$defWs.Range("C2").Value = 156
$defWs.Range("D2").Value = 14
$defWs.Range("F2").Value = 55
$defWs.Range("G2").Value = 38
$defWs.Range("H2").Value = 4
$defWs.Range("J2").Value = 25
$defWs.Range("L2").Value = 248
$defWs.Range("M2").Value = 160
$defWs.Range("O2").Value = 21
$defWs.Range("P2").Value = 13
$defWs.Range("Q2").Value = 472

$defWs.Range("B3").Value = 6
$defWs.Range("C3").Value = 156
$defWs.Range("E3").Value = 20
$defWs.Range("F3").Value = 102
$defWs.Range("G3").Value = 31
$defWs.Range("H3").Value = 30
$defWs.Range("I3").Value = 58
$defWs.Range("J3").Value = 36
$defWs.Range("N3").Value = 22

# ---------------------------------------------------------------------
# ST sheet: append Week 16 logs + updated season totals
# ---------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B2").Value = 62
$stWs.Range("D2").Value = 50
$stWs.Range("F2").Value = 64
$stWs.Range("G2").Value = 62
$stWs.Range("L2").Value = 21

$stWs.Range("B3").Value = 27

Append-Numbers $stWs "B4" "52 63"
Append-Numbers $stWs "B5" "14 23"
Append-Numbers $stWs "B6" "24"
Append-Numbers $stWs "D3" "49 51 44 37 43"
Append-Numbers $stWs "D4" "29 0 0 0 2"
Append-Numbers $stWs "D5" "20 0 0 0"

# ---------------------------------------------------------------------
# TURNS sheet: updated season totals
# ---------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")

$turnsWs.Range("C2").Value = 6
$turnsWs.Range("E2").Value = 10
$turnsWs.Range("E3").Value = 7

# ---------------------------------------------------------------------
# PEN sheet: updated season totals
# ---------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")

$penWs.Range("B3").Value = 24
$penWs.Range("D4").Value = 11

$wb.Save()
